# Weekly fruit/vegetable price update: insert a new weekly record as row 7
# (Fecha 2023-10-31 / serial 45230), pushing the existing rows 7-21 down to
# rows 8-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 7 - this shifts rows 7..21 down
# to 8..22 and copies the row's formatting (e.g. the date-number style on
# column D) from the row above, just like an interactive Excel insert.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly entry.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 45230
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 300000000
$ws.Cells.Item(7, 7).Value = "Espárragos"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1500
$ws.Cells.Item(7, 14).Value = "`$/kilo"
$ws.Cells.Item(7, 15).Value = "Provincia de Linares"
$ws.Cells.Item(7, 16).Value = 1500
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
